# Automatische test-sync: 2025-06-26 21:21:50
# Appends the new mail-log entry (row 18) to the "Logs" sheet and bumps the
# "Bestelling / Levering" count on the "Dashboard" sheet from 12 to 13.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$row = 18

$logs.Cells.Item($row, 1).Value = "Wil je 5 boren bestellen?"
$logs.Cells.Item($row, 2).Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Cells.Item($row, 3).Value = "Hoi Johan,`nWil je 5 boren bestellen?`nRick`nSent using {0}"
$logs.Cells.Item($row, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($row, 5).Value = "Beste Rick,`nBedankt voor je e-mail. Om de bestelling van de 5 boren te verwerken, hebben we wat meer informatie nodig. Zou je het artikelnummer of de specifieke kenmerken van de boren kunnen doorgeven? Op die manier kunnen we ervoor zorgen dat we de juiste producten voor je bestellen.`nMet vriendelijke groet,`nJohan"
$logs.Cells.Item($row, 6).Value = "2025-06-26 21:21:18"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Ja"

# Entering multi-line text (hard line breaks) bumps this row's height away
# from the sheet default; re-running AutoFit snaps it back to a plain,
# un-pinned row (no explicit ht/customHeight), matching the source rows.
$logs.Rows.Item($row).AutoFit()

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 13

# Extend the conditional-formatting ranges (D/G/H/I 2:17 -> 2:18) to cover
# the newly appended row.
$logs.Range("D2:D17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D18"))
$logs.Range("G2:G17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G18"))
$logs.Range("H2:H17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H18"))
$logs.Range("I2:I17").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I18"))
